# This script updates the "想去人数" (number of people interested) column (F)
# on the "展览" (sheet index 1) and "全部类型" (sheet index 4) worksheets,
# reflecting a refreshed data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Update "想去人数" (interested count) values in sheet "展览" (index 1)
$ws1.Cells.Item(5, 6).Value = 12
$ws1.Cells.Item(6, 6).Value = 200
$ws1.Cells.Item(7, 6).Value = 4611
$ws1.Cells.Item(9, 6).Value = 122
$ws1.Cells.Item(12, 6).Value = 90
$ws1.Cells.Item(13, 6).Value = 707
$ws1.Cells.Item(14, 6).Value = 191
$ws1.Cells.Item(15, 6).Value = 1010
$ws1.Cells.Item(16, 6).Value = 80
$ws1.Cells.Item(17, 6).Value = 242
$ws1.Cells.Item(18, 6).Value = 157
$ws1.Cells.Item(19, 6).Value = 72
$ws1.Cells.Item(20, 6).Value = 121
$ws1.Cells.Item(22, 6).Value = 3581
$ws1.Cells.Item(23, 6).Value = 5930
$ws1.Cells.Item(27, 6).Value = 524
$ws1.Cells.Item(29, 6).Value = 3375
$ws1.Cells.Item(30, 6).Value = 368
$ws1.Cells.Item(32, 6).Value = 2480
$ws1.Cells.Item(36, 6).Value = 218
$ws1.Cells.Item(37, 6).Value = 270
$ws1.Cells.Item(38, 6).Value = 354
$ws1.Cells.Item(39, 6).Value = 132
$ws1.Cells.Item(40, 6).Value = 1506
$ws1.Cells.Item(41, 6).Value = 915
$ws1.Cells.Item(43, 6).Value = 30
$ws1.Cells.Item(44, 6).Value = 49
$ws1.Cells.Item(46, 6).Value = 472
$ws1.Cells.Item(47, 6).Value = 66
$ws1.Cells.Item(48, 6).Value = 553

# Update "想去人数" (interested count) values in sheet "全部类型" (index 4)
$ws4.Cells.Item(5, 6).Value = 12
$ws4.Cells.Item(6, 6).Value = 200
$ws4.Cells.Item(7, 6).Value = 4611
$ws4.Cells.Item(9, 6).Value = 122
$ws4.Cells.Item(13, 6).Value = 90
$ws4.Cells.Item(14, 6).Value = 707
$ws4.Cells.Item(15, 6).Value = 191
$ws4.Cells.Item(16, 6).Value = 1010
$ws4.Cells.Item(17, 6).Value = 80
$ws4.Cells.Item(18, 6).Value = 242
$ws4.Cells.Item(19, 6).Value = 157
$ws4.Cells.Item(20, 6).Value = 72
$ws4.Cells.Item(21, 6).Value = 121
$ws4.Cells.Item(23, 6).Value = 3581
$ws4.Cells.Item(24, 6).Value = 5930
$ws4.Cells.Item(28, 6).Value = 524
$ws4.Cells.Item(30, 6).Value = 3375
$ws4.Cells.Item(31, 6).Value = 368
$ws4.Cells.Item(33, 6).Value = 2480
$ws4.Cells.Item(37, 6).Value = 218
$ws4.Cells.Item(38, 6).Value = 270
$ws4.Cells.Item(39, 6).Value = 354
$ws4.Cells.Item(40, 6).Value = 132
$ws4.Cells.Item(41, 6).Value = 1508
$ws4.Cells.Item(42, 6).Value = 915
$ws4.Cells.Item(44, 6).Value = 30
$ws4.Cells.Item(45, 6).Value = 49
$ws4.Cells.Item(47, 6).Value = 472
$ws4.Cells.Item(48, 6).Value = 66
$ws4.Cells.Item(49, 6).Value = 553
